# Create main thread table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns A:C to be stored as text so that numeric-looking and
# date-looking values (e.g. "37", "2019-08-07 15:36:04.947") are kept
# as literal strings instead of being converted to numbers/dates.
$ws.Range("A1:C11").NumberFormat = "@"

$data = @(
    @("主线程开始", "", "2019-08-07 15:36:04.947"),
    @("basicOpe", "37", "2019-08-07 15:36:04.998"),
    @("basicProd", "3", "2019-08-07 15:36:05.001"),
    @("basicOwnTime", "15", "2019-08-07 15:36:05.016"),
    @("rpushTime", "3", "2019-08-07 15:36:05.019"),
    @("zaddTime", "1", "2019-08-07 15:36:05.020"),
    @("matchedRuleTime", "158", "2019-08-07 15:36:05.178"),
    @("zaddHistory", "166", "2019-08-07 15:36:05.344"),
    @("主线程调用子线程(MAIN CALL SUB)", "", "2019-08-07 15:36:05.496"),
    @("回主线程 (MAIN RETURN)", "", "2019-08-07 15:36:09.602"),
    @("总时间", "4106", "")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
}
